$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2363.3333
$ws.Range("I6").Value = 135
$ws.Range("K6").Value = 405
$ws.Range("M6").Value = -293
$ws.Range("H49").Value = 50
$ws.Range("J49").Value = 50
$ws.Range("L49").Value = 150
$ws.Range("N49").Value = -422
$ws.Range("H87").Value = 39982
$ws.Range("J87").Value = 39982
$ws.Range("L87").Value = 39982
$ws.Range("N87").Value = -42478
$ws.Range("H90").Value = 39982
$ws.Range("J90").Value = 39982
$ws.Range("L90").Value = 119946
$ws.Range("N90").Value = -132426
$ws.Range("H106").Value = 6242.25
$ws.Range("I106").Value = 6323
$ws.Range("K106").Value = 6323
$ws.Range("M106").Value = -5692
$ws.Range("H137").Value = 3749.5
$ws.Range("J137").Value = 2500
$ws.Range("L137").Value = 7500
$ws.Range("N137").Value = -12600
$ws.Range("H138").Value = 8193.5
$ws.Range("J138").Value = 8219.733
$ws.Range("L138").Value = 24659.199
$ws.Range("N138").Value = -34939.199

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4078.2917
$ws.Range("I32").Value = 3585.9092
$ws.Range("K32").Value = 3585.9092
$ws.Range("M32").Value = -3298.9092
$ws.Range("H53").Value = 10000
$ws.Range("I53").Value = 10000
$ws.Range("K53").Value = 10000
$ws.Range("M53").Value = -9318
$ws.Range("H122").Value = 1782.6818
$ws.Range("I122").Value = 1782.6818
$ws.Range("K122").Value = 5348.0454
$ws.Range("M122").Value = -2898.0454

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 20000
$ws.Range("J8").Value = 20000
$ws.Range("L8").Value = 20000
$ws.Range("N8").Value = -20280
$ws.Range("H86").Value = 2475
$ws.Range("I86").Value = 2300
$ws.Range("J86").Value = 2650
$ws.Range("K86").Value = 2300
$ws.Range("L86").Value = 2650
$ws.Range("M86").Value = -1177
$ws.Range("N86").Value = -4896
$ws.Range("H89").Value = 2475
$ws.Range("I89").Value = 2300
$ws.Range("J89").Value = 2650
$ws.Range("K89").Value = 11500
$ws.Range("L89").Value = 13250
$ws.Range("M89").Value = -5884
$ws.Range("N89").Value = -24482
$ws.Range("H134").Value = 1910.8823
$ws.Range("I134").Value = 1905.8
$ws.Range("K134").Value = 5717.4
$ws.Range("M134").Value = -3182.4

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 50092
$ws.Range("J50").Value = 50092
$ws.Range("L50").Value = 50092
$ws.Range("N50").Value = -51342
$ws.Range("H51").Value = 40049.5
$ws.Range("J51").Value = 50099
$ws.Range("L51").Value = 50099
$ws.Range("N51").Value = -51571
$ws.Range("H60").Value = 36594.715
$ws.Range("J60").Value = 36594.715
$ws.Range("L60").Value = 36594.715
$ws.Range("N60").Value = -37616.715
$ws.Range("H61").Value = 40049.5
$ws.Range("J61").Value = 50099
$ws.Range("L61").Value = 50099
$ws.Range("N61").Value = -50795
$ws.Range("H132").Value = 2831.3
$ws.Range("I132").Value = 1720.0834
$ws.Range("K132").Value = 5160.2502
$ws.Range("M132").Value = -2630.2502
$ws.Range("H134").Value = 1677.6875
$ws.Range("I134").Value = 1346
$ws.Range("J134").Value = 3999.5
$ws.Range("K134").Value = 4038
$ws.Range("L134").Value = 11998.5
$ws.Range("M134").Value = -1503
$ws.Range("N134").Value = -17068.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 76.2
$ws.Range("I2").Value = 57.666668
$ws.Range("J2").Value = 104
$ws.Range("K2").Value = 346.000008
$ws.Range("L2").Value = 624
$ws.Range("M2").Value = -233.000008
$ws.Range("N2").Value = -850
$ws.Range("H4").Value = 137645390
$ws.Range("I4").Value = 105078780
$ws.Range("K4").Value = 315236340
$ws.Range("M4").Value = -315236228
$ws.Range("H47").Value = 649.5
$ws.Range("I47").Value = 649.5
$ws.Range("K47").Value = 1948.5
$ws.Range("M47").Value = -1517.5
$ws.Range("H92").Value = 1167
$ws.Range("I92").Value = 1167
$ws.Range("K92").Value = 3501
$ws.Range("M92").Value = -2253
$ws.Range("H94").Value = 2899
$ws.Range("J94").Value = 2899
$ws.Range("L94").Value = 8697
$ws.Range("N94").Value = -10049
$ws.Range("H107").Value = 4619.8
$ws.Range("J107").Value = 1033
$ws.Range("L107").Value = 3099
$ws.Range("N107").Value = -6939
$ws.Range("H109").Value = 998.6667
$ws.Range("J109").Value = 998
$ws.Range("L109").Value = 2994
$ws.Range("N109").Value = -5074
$ws.Range("H111").Value = 824
$ws.Range("I111").Value = 824
$ws.Range("K111").Value = 2472
$ws.Range("M111").Value = 595

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 11799.5
$ws.Range("J96").Value = 11799.5
$ws.Range("L96").Value = 11799.5
$ws.Range("N96").Value = -17291.5
$ws.Range("H102").Value = 4225.143
$ws.Range("I102").Value = 3262.6667
$ws.Range("K102").Value = 3262.6667
$ws.Range("M102").Value = -1640.6667
$ws.Range("H122").Value = 1999.75
$ws.Range("I122").Value = 1999.75
$ws.Range("K122").Value = 5999.25
$ws.Range("M122").Value = -3549.25
$ws.Range("H132").Value = 4109.706
$ws.Range("I132").Value = 3704.9285
$ws.Range("J132").Value = 5998.6665
$ws.Range("K132").Value = 11114.7855
$ws.Range("L132").Value = 17995.9995
$ws.Range("M132").Value = -8584.7855
$ws.Range("N132").Value = -23055.9995

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3659.375
$ws.Range("J40").Value = 3594.5
$ws.Range("L40").Value = 3594.5
$ws.Range("N40").Value = -3866.5
$ws.Range("H132").Value = 3041.7368
$ws.Range("I132").Value = 2446.2307
$ws.Range("K132").Value = 7338.6921
$ws.Range("M132").Value = -4808.6921
$ws.Range("H136").Value = 29411.125
$ws.Range("I136").Value = 5989.4546
$ws.Range("K136").Value = 17968.3638
$ws.Range("M136").Value = -15418.3638

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2690.7273
$ws.Range("I132").Value = 2271.1667
$ws.Range("K132").Value = 6813.500100000001
$ws.Range("M132").Value = -4283.500100000001
